$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "earn_rule_outcome_1565160"
$ws.Range("B3").Select()
